$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 9300
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -10224

$ws.Range("H9").Value = 5561513.5
$ws.Range("I9").Value = 11905278
$ws.Range("K9").Value = 11905278
$ws.Range("M9").Value = -11905109

$ws.Range("H14").Value = 9300
$ws.Range("I14").Value = 3000
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -2809
$ws.Range("N14").Value = -10382

$ws.Range("H17").Value = 130711.49
$ws.Range("J17").Value = 133305.72
$ws.Range("L17").Value = 399917.16
$ws.Range("N17").Value = -400253.16

$ws.Range("H18").Value = 3924.7273
$ws.Range("I18").Value = 3574.6667
$ws.Range("K18").Value = 3574.6667
$ws.Range("M18").Value = -3290.6667

$ws.Range("H40").Value = 1952.3077
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 1988
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 1988
$ws.Range("M40").Value = -1658.3334
$ws.Range("N40").Value = -2338

$ws.Range("H63").Value = 28000
$ws.Range("J63").Value = 28000
$ws.Range("L63").Value = 28000
$ws.Range("N63").Value = -29248

$ws.Range("H66").Value = 28000
$ws.Range("J66").Value = 28000
$ws.Range("L66").Value = 84000
$ws.Range("N66").Value = -90240

$ws.Range("H86").Value = 47622344
$ws.Range("I86").Value = 62503460
$ws.Range("J86").Value = 2779.4
$ws.Range("K86").Value = 62503460
$ws.Range("L86").Value = 2779.4
$ws.Range("M86").Value = -62502337
$ws.Range("N86").Value = -5025.4

$ws.Range("H88").Value = 3412.8667
$ws.Range("J88").Value = 4156.1816
$ws.Range("L88").Value = 4156.1816
$ws.Range("N88").Value = -4968.1816

$ws.Range("H89").Value = 47622344
$ws.Range("I89").Value = 62503460
$ws.Range("J89").Value = 2779.4
$ws.Range("K89").Value = 312517300
$ws.Range("L89").Value = 13897
$ws.Range("M89").Value = -312511684
$ws.Range("N89").Value = -25129

$ws.Range("H91").Value = 3412.8667
$ws.Range("J91").Value = 4156.1816
$ws.Range("L91").Value = 4156.1816
$ws.Range("N91").Value = -6964.1816

$ws.Range("H113").Value = 2792.6365
$ws.Range("I113").Value = 2341.5
$ws.Range("J113").Value = 3995.6667
$ws.Range("K113").Value = 2341.5
$ws.Range("L113").Value = 3995.6667
$ws.Range("M113").Value = 912.5
$ws.Range("N113").Value = -10503.6667

$ws.Range("H116").Value = 14617
$ws.Range("J116").Value = 4198
$ws.Range("L116").Value = 4198
$ws.Range("N116").Value = -11082

$ws.Range("H132").Value = 42984.668
$ws.Range("I132").Value = 46710.547
$ws.Range("K132").Value = 140131.641
$ws.Range("M132").Value = -137601.641

$ws.Range("H137").Value = 2500974.8
$ws.Range("I137").Value = 977.2258
$ws.Range("J137").Value = 11112078
$ws.Range("K137").Value = 2931.6774
$ws.Range("L137").Value = 33336234
$ws.Range("M137").Value = -381.6774
$ws.Range("N137").Value = -33341334

$ws.Range("H138").Value = 5135.164
$ws.Range("I138").Value = 11208.857
$ws.Range("J138").Value = 3325.9788
$ws.Range("K138").Value = 33626.571
$ws.Range("L138").Value = 9977.936399999999
$ws.Range("M138").Value = -28486.571
$ws.Range("N138").Value = -20257.9364

$ws.Range("H141").Value = 1186.4445
$ws.Range("I141").Value = 978.5
$ws.Range("K141").Value = 2935.5
$ws.Range("M141").Value = 2244.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 8750
$ws.Range("I10").Value = 7500
$ws.Range("K10").Value = 7500
$ws.Range("M10").Value = -7330

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()

$ws.Range("H45").Value = 62604.168
$ws.Range("I45").Value = 116875
$ws.Range("K45").Value = 116875
$ws.Range("M45").Value = -116498

$ws.Range("H74").Value = 339937.72
$ws.Range("I74").Value = 1659.0769
$ws.Range("J74").Value = 828562.4399999999
$ws.Range("K74").Value = 1659.0769
$ws.Range("L74").Value = 828562.4399999999
$ws.Range("M74").Value = -785.0769
$ws.Range("N74").Value = -830310.4399999999

$ws.Range("H77").Value = 339937.72
$ws.Range("I77").Value = 1659.0769
$ws.Range("J77").Value = 828562.4399999999
$ws.Range("K77").Value = 8295.3845
$ws.Range("L77").Value = 4142812.2
$ws.Range("M77").Value = -3927.3845
$ws.Range("N77").Value = -4151548.2

$ws.Range("H88").Value = 1766.0714
$ws.Range("I88").Value = 1820.25
$ws.Range("J88").Value = 1744.4
$ws.Range("K88").Value = 1820.25
$ws.Range("L88").Value = 1744.4
$ws.Range("M88").Value = -1414.25
$ws.Range("N88").Value = -2556.4

$ws.Range("H91").Value = 1766.0714
$ws.Range("I91").Value = 1820.25
$ws.Range("J91").Value = 1744.4
$ws.Range("K91").Value = 1820.25
$ws.Range("L91").Value = 1744.4
$ws.Range("M91").Value = -416.25
$ws.Range("N91").Value = -4552.4

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H96").Value = 46420
$ws.Range("J96").Value = 46420
$ws.Range("L96").Value = 46420
$ws.Range("N96").Value = -51912

$ws.Range("H132").Value = 2540.7144
$ws.Range("J132").Value = 3530.1667
$ws.Range("L132").Value = 10590.5001
$ws.Range("N132").Value = -15650.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7262.65
$ws.Range("I105").Value = 7777.067
$ws.Range("K105").Value = 7777.067
$ws.Range("M105").Value = -6030.067

$ws.Range("H107").Value = 9004.833000000001
$ws.Range("I107").Value = 11194.546
$ws.Range("K107").Value = 11194.546
$ws.Range("M107").Value = -9274.546

$ws.Range("H134").Value = 25716898
$ws.Range("I134").Value = 2450.238
$ws.Range("J134").Value = 64288572
$ws.Range("K134").Value = 7350.714
$ws.Range("L134").Value = 192865716
$ws.Range("M134").Value = -4815.714
$ws.Range("N134").Value = -192870786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2431.35
$ws.Range("I31").Value = 1696.2759
$ws.Range("K31").Value = 1696.2759
$ws.Range("M31").Value = -1401.2759

$ws.Range("H34").Value = 2431.35
$ws.Range("I34").Value = 1696.2759
$ws.Range("K34").Value = 1696.2759
$ws.Range("M34").Value = -1494.2759

$ws.Range("H43").Value = 12290.637
$ws.Range("J43").Value = 12290.637
$ws.Range("L43").Value = 12290.637
$ws.Range("N43").Value = -12658.637

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H101").Value = 12290.637
$ws.Range("J101").Value = 12290.637
$ws.Range("L101").Value = 12290.637
$ws.Range("N101").Value = -18780.637

$ws.Range("H105").Value = 1922.9286
$ws.Range("I105").Value = 1447.5454
$ws.Range("K105").Value = 1447.5454
$ws.Range("M105").Value = 299.4546

$ws.Range("H107").Value = 1357.7667
$ws.Range("I107").Value = 1400.9546
$ws.Range("J107").Value = 1239
$ws.Range("K107").Value = 1400.9546
$ws.Range("L107").Value = 1239
$ws.Range("M107").Value = 519.0454
$ws.Range("N107").Value = -5079

$ws.Range("H132").Value = 22016.06
$ws.Range("I132").Value = 30141.6
$ws.Range("J132").Value = 1702.2142
$ws.Range("K132").Value = 90424.79999999999
$ws.Range("L132").Value = 5106.642599999999
$ws.Range("M132").Value = -87894.79999999999
$ws.Range("N132").Value = -10166.6426

$ws.Range("H134").Value = 2285.739
$ws.Range("I134").Value = 1750.6666
$ws.Range("K134").Value = 5251.9998
$ws.Range("M134").Value = -2716.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 496.6154
$ws.Range("J38").Value = 701.3333
$ws.Range("L38").Value = 2103.9999
$ws.Range("N38").Value = -2797.9999

$ws.Range("H92").Value = 353.875
$ws.Range("I92").Value = 241.5
$ws.Range("K92").Value = 724.5
$ws.Range("M92").Value = 523.5

$ws.Range("H97").Value = 290.7143
$ws.Range("I97").Value = 107.2
$ws.Range("K97").Value = 321.6
$ws.Range("M97").Value = 174.4

$ws.Range("H101").Value = 20995
$ws.Range("J101").Value = 20995
$ws.Range("L101").Value = 62985
$ws.Range("N101").Value = -67853

$ws.Range("H109").Value = 7357.8
$ws.Range("I109").Value = 4197.25
$ws.Range("K109").Value = 12591.75
$ws.Range("M109").Value = -11551.75

$ws.Range("H134").Value = 4506.385
$ws.Range("I134").Value = 2598.4546
$ws.Range("K134").Value = 7795.3638
$ws.Range("M134").Value = -2725.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15626865
$ws.Range("I102").Value = 20001668
$ws.Range("K102").Value = 20001668
$ws.Range("M102").Value = -20000046

$ws.Range("H113").Value = 1968.7333
$ws.Range("I113").Value = 1714.5555
$ws.Range("K113").Value = 1714.5555
$ws.Range("M113").Value = 455.4445000000001

$ws.Range("H122").Value = 1652.8605
$ws.Range("I122").Value = 1644.7241
$ws.Range("J122").Value = 1669.7142
$ws.Range("K122").Value = 4934.1723
$ws.Range("L122").Value = 5009.142599999999
$ws.Range("M122").Value = -2484.1723
$ws.Range("N122").Value = -9909.142599999999

$ws.Range("H132").Value = 390005.97
$ws.Range("I132").Value = 1212.2424
$ws.Range("K132").Value = 3636.7272
$ws.Range("M132").Value = -1106.7272

$ws.Range("H137").Value = 176626
$ws.Range("I137").Value = 169898
$ws.Range("K137").Value = 169898
$ws.Range("M137").Value = -164798

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1112.2593
$ws.Range("I16").Value = 1075.6316
$ws.Range("J16").Value = 1199.25
$ws.Range("K16").Value = 1075.6316
$ws.Range("L16").Value = 1199.25
$ws.Range("M16").Value = -905.6315999999999
$ws.Range("N16").Value = -1539.25

$ws.Range("H40").Value = 1710.1177
$ws.Range("I40").Value = 961.63635
$ws.Range("K40").Value = 961.63635
$ws.Range("M40").Value = -825.63635

$ws.Range("H46").Value = 10159
$ws.Range("I46").Value = 13085.454
$ws.Range("K46").Value = 13085.454
$ws.Range("M46").Value = -12897.454

$ws.Range("H55").Value = 1118.9565
$ws.Range("I55").Value = 1061.5
$ws.Range("K55").Value = 1061.5
$ws.Range("M55").Value = -888.5

$ws.Range("H93").Value = 1180.7273
$ws.Range("I93").Value = 1161.9474
$ws.Range("J93").Value = 1299.6666
$ws.Range("K93").Value = 1161.9474
$ws.Range("L93").Value = 1299.6666
$ws.Range("M93").Value = 86.05259999999998
$ws.Range("N93").Value = -3795.6666

$ws.Range("H122").Value = 2163.5908
$ws.Range("I122").Value = 1881.6666
$ws.Range("K122").Value = 5644.9998
$ws.Range("M122").Value = -3194.9998

$ws.Range("H132").Value = 2236.3394
$ws.Range("I132").Value = 1731.3684
$ws.Range("K132").Value = 5194.1052
$ws.Range("M132").Value = -2664.1052

$ws.Range("H137").Value = 119995
$ws.Range("J137").Value = 189990
$ws.Range("L137").Value = 189990
$ws.Range("N137").Value = -200190

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H113").Value = 543.43475
$ws.Range("I113").Value = 434.52942
$ws.Range("K113").Value = 1303.58826
$ws.Range("M113").Value = 866.41174

$ws.Range("H126").Value = 2140.7778
$ws.Range("I126").Value = 1871.6666
$ws.Range("K126").Value = 5614.9998
$ws.Range("M126").Value = -3144.9998

$ws.Range("H132").Value = 2156.8113
$ws.Range("J132").Value = 2343.4167
$ws.Range("L132").Value = 7030.250100000001
$ws.Range("N132").Value = -12090.2501

